$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Calificaciones Grupo 1")
$ws2 = $wb.Worksheets.Item("Calificaciones Gr 51")

# --- Update "tutorias" (AB column) scores; dependent formulas (AC/AD/AE/AH)
# recalculate automatically. ---
$ws1.Range("AB7").Value = 10
$ws1.Range("AB13").Value = 6
$ws1.Range("AB14").Value = 10
$ws1.Range("AB15").Value = 4
$ws1.Range("AB16").Value = 10
$ws1.Range("AB17").Value = 2
$ws1.Range("AB24").Value = 10

$ws2.Range("AB3").Value = 10

# --- T4 on sheet 2 picks up the same (visually identical) cell style used
# by other cells in column T, e.g. T14, matching the author's re-entry of
# the cell during editing. ---
$ws2.Range("T14").Copy()
$ws2.Range("T4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View/selection state: the author ended the session on sheet 1 with a
# different scroll position and selection than before; sheet 2's zoom got
# reset back to 100% (default) and its selection/scroll moved too. ---
$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 20
$excel.ActiveWindow.Zoom = 100
$ws2.Range("AH20").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
$ws1.Range("AE7").Select()
